$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3001
$ws.Range("J32").Value = 5002
$ws.Range("L32").Value = 5002
$ws.Range("N32").Value = -5654

$ws.Range("H62").Value = 5211.222
$ws.Range("I62").Value = 1399.1666
$ws.Range("J62").Value = 12835.333
$ws.Range("K62").Value = 1399.1666
$ws.Range("L62").Value = 12835.333
$ws.Range("M62").Value = -775.1666
$ws.Range("N62").Value = -14083.333

$ws.Range("H65").Value = 5211.222
$ws.Range("I65").Value = 1399.1666
$ws.Range("J65").Value = 12835.333
$ws.Range("K65").Value = 6995.833000000001
$ws.Range("L65").Value = 64176.665
$ws.Range("M65").Value = -3875.833000000001
$ws.Range("N65").Value = -70416.66500000001

$ws.Range("H137").Value = 26317056
$ws.Range("J137").Value = 2864.1667
$ws.Range("L137").Value = 8592.500100000001
$ws.Range("N137").Value = -13692.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18356.436
$ws.Range("I32").Value = 2373.3035
$ws.Range("K32").Value = 2373.3035
$ws.Range("M32").Value = -2086.3035

$ws.Range("H61").Value = 3203.5
$ws.Range("J61").Value = 3997.8
$ws.Range("L61").Value = 3997.8
$ws.Range("N61").Value = -4421.8

$ws.Range("H122").Value = 1602.3143
$ws.Range("I122").Value = 1240.7142
$ws.Range("J122").Value = 2144.7144
$ws.Range("K122").Value = 3722.1426
$ws.Range("L122").Value = 6434.1432
$ws.Range("M122").Value = -1272.1426
$ws.Range("N122").Value = -11334.1432

$ws.Range("H132").Value = 2735.5166
$ws.Range("I132").Value = 2494.577
$ws.Range("J132").Value = 4301.625
$ws.Range("K132").Value = 7483.731000000001
$ws.Range("L132").Value = 12904.875
$ws.Range("M132").Value = -4953.731000000001
$ws.Range("N132").Value = -17964.875

$ws.Range("H136").Value = 3203.5
$ws.Range("J136").Value = 3997.8
$ws.Range("L136").Value = 11993.4
$ws.Range("N136").Value = -17093.4

$ws.Range("H138").Value = 55000
$ws.Range("J138").Value = 55000
$ws.Range("L138").Value = 55000
$ws.Range("N138").Value = -65280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2240.682
$ws.Range("J20").Value = 1790
$ws.Range("L20").Value = 1790
$ws.Range("N20").Value = -2284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 31857.715
$ws.Range("I62").Value = 42601
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 42601
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -41977
$ws.Range("N62").Value = -6247.5

$ws.Range("H65").Value = 31857.715
$ws.Range("I65").Value = 42601
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 213005
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -209885
$ws.Range("N65").Value = -31237.5

$ws.Range("H132").Value = 3431.5217
$ws.Range("I132").Value = 2958.2856
$ws.Range("J132").Value = 4167.6665
$ws.Range("K132").Value = 8874.856800000001
$ws.Range("L132").Value = 12502.9995
$ws.Range("M132").Value = -6344.856800000001
$ws.Range("N132").Value = -17562.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3139.8
$ws.Range("I70").Value = 1899.6666
$ws.Range("K70").Value = 5698.9998
$ws.Range("M70").Value = -5383.9998

$ws.Range("H73").Value = 3139.8
$ws.Range("I73").Value = 1899.6666
$ws.Range("K73").Value = 5698.9998
$ws.Range("M73").Value = -4606.9998

$ws.Range("H75").Value = 2550.3333
$ws.Range("I75").Value = 1575.5
$ws.Range("J75").Value = 4500
$ws.Range("K75").Value = 4726.5
$ws.Range("L75").Value = 13500
$ws.Range("M75").Value = -3728.5
$ws.Range("N75").Value = -15496

$ws.Range("H78").Value = 2550.3333
$ws.Range("I78").Value = 1575.5
$ws.Range("J78").Value = 4500
$ws.Range("K78").Value = 14179.5
$ws.Range("L78").Value = 40500
$ws.Range("M78").Value = -9187.5
$ws.Range("N78").Value = -50484

$ws.Range("H103").Value = 1807
$ws.Range("I103").Value = 1149.375
$ws.Range("J103").Value = 2211.6924
$ws.Range("K103").Value = 3448.125
$ws.Range("L103").Value = 6635.0772
$ws.Range("M103").Value = -2569.125
$ws.Range("N103").Value = -8393.0772

$ws.Range("H113").Value = 907.9355
$ws.Range("I113").Value = 633.44446
$ws.Range("J113").Value = 1020.2273
$ws.Range("K113").Value = 1900.33338
$ws.Range("L113").Value = 3060.6819
$ws.Range("M113").Value = 269.66662
$ws.Range("N113").Value = -7400.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""

$ws.Range("H102").Value = 1900.2174
$ws.Range("I102").Value = 1738.5
$ws.Range("K102").Value = 1738.5
$ws.Range("M102").Value = -116.5

$ws.Range("H107").Value = 600.44446
$ws.Range("I107").Value = 557.2857
$ws.Range("J107").Value = 751.5
$ws.Range("K107").Value = 557.2857
$ws.Range("L107").Value = 751.5
$ws.Range("M107").Value = 1362.7143
$ws.Range("N107").Value = -4591.5

$ws.Range("H116").Value = 30000
$ws.Range("J116").Value = 30000
$ws.Range("L116").Value = 30000
$ws.Range("N116").Value = -39178

$ws.Range("H132").Value = 3953.3845
$ws.Range("I132").Value = 3781.7144
$ws.Range("K132").Value = 11345.1432
$ws.Range("M132").Value = -8815.143199999999

$ws.Range("H140").Value = 55375
$ws.Range("J140").Value = 55375
$ws.Range("L140").Value = 55375
$ws.Range("N140").Value = -65735

$ws.Range("H141").Value = 62742.832
$ws.Range("J141").Value = 62742.832
$ws.Range("L141").Value = 62742.832
$ws.Range("N141").Value = -73102.832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1041.6154
$ws.Range("I46").Value = 1029.625
$ws.Range("J46").Value = 1060.8
$ws.Range("K46").Value = 1029.625
$ws.Range("L46").Value = 1060.8
$ws.Range("M46").Value = -841.625
$ws.Range("N46").Value = -1436.8

$ws.Range("H93").Value = 1249.2307
$ws.Range("I93").Value = 1117.75
$ws.Range("J93").Value = 1459.6
$ws.Range("K93").Value = 1117.75
$ws.Range("L93").Value = 1459.6
$ws.Range("M93").Value = 130.25
$ws.Range("N93").Value = -3955.6

$ws.Range("H122").Value = 3528.8262
$ws.Range("I122").Value = 1802
$ws.Range("J122").Value = 3892.3684
$ws.Range("K122").Value = 5406
$ws.Range("L122").Value = 11677.1052
$ws.Range("M122").Value = -2956
$ws.Range("N122").Value = -16577.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 299.8095
$ws.Range("I100").Value = 206.46666
$ws.Range("J100").Value = 533.1667
$ws.Range("K100").Value = 412.93332
$ws.Range("L100").Value = 1066.3334
$ws.Range("M100").Value = 128.06668
$ws.Range("N100").Value = -2148.3334

$ws.Range("H122").Value = 32703.94
$ws.Range("I122").Value = 45274.043
$ws.Range("J122").Value = 3792.7
$ws.Range("K122").Value = 135822.129
$ws.Range("L122").Value = 11378.1
$ws.Range("M122").Value = -133372.129
$ws.Range("N122").Value = -16278.1

$ws.Range("H141").Value = 75357.5
$ws.Range("J141").Value = 75357.5
$ws.Range("L141").Value = 75357.5
$ws.Range("N141").Value = -85717.5

Write-Output "Applied all changes"
